$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.160.79'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.03%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.255.76'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.53%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '99.00'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +17.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '272.14'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.66%  '

$ws.Range('E7').Value = '  +0.94%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.627'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.92%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.94'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +6.99%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0946'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.35%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.20'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +13.85%  '

$ws.Range('E13').Value = '  +0.32%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.47'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +7.70%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.590.08'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.29%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.837'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +7.21%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.242.82'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.64%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.114.47'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.11%  '

$ws.Range('E19').Value = '  +3.89%  '

$ws.Range('E20').Value = '  +5.13%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.82'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.39%  '

$ws.Range('E22').Value = '  +0.92%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.04'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +9.69%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '234.69'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.29%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.38'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.70%  '

$ws.Range('E27').Value = '  +12.09%  '

$ws.Range('E28').Value = '  -3.05%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.01'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.03%  '

$ws.Range('E30').Value = '  +0.66%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.58'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.07%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0914'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.46%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.23'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.80%  '

$ws.Range('E35').Value = '  +1.64%  '

$ws.Range('E36').Value = '  +0.90%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0356'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.02%  '

$ws.Range('E38').Value = '  -3.25%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.54'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +24.45%  '

$ws.Range('E40').Value = '  +24.34%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.19'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.42%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '12.54'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.55%  '

$ws.Range('E43').Value = '  -0.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.10'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.43%  '

$ws.Range('E45').Value = '  +5.24%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.53'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.50%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.54'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.37%  '

$ws.Range('E48').Value = '  +4.36%  '

$ws.Range('E49').Value = '  -0.15%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.430'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.08%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.471.57'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.16%  '
